# Updates cryptos price/volume figures per the commit diff.
# Uses a helper cell + text-formula + PasteSpecial(xlPasteValues) trick so that
# numeric-looking strings (e.g. "64.30") are written as literal text, matching
# the workbook's existing convention of storing these columns as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '58.156.35' },
    @{ Cell = 'D3'; Value = '2.509.52' },
    @{ Cell = 'E3'; Value = '  +1.93%  ' },
    @{ Cell = 'D4'; Value = '0.999' },
    @{ Cell = 'E4'; Value = '  -0.16%  ' },
    @{ Cell = 'D5'; Value = '521.32' },
    @{ Cell = 'E5'; Value = '  +0.21%  ' },
    @{ Cell = 'D6'; Value = '131.95' },
    @{ Cell = 'E6'; Value = '  -0.86%  ' },
    @{ Cell = 'D7'; Value = '0.997' },
    @{ Cell = 'E7'; Value = '  -0.23%  ' },
    @{ Cell = 'D8'; Value = '0.556' },
    @{ Cell = 'E8'; Value = '  +0.04%  ' },
    @{ Cell = 'D9'; Value = '2.508.57' },
    @{ Cell = 'E9'; Value = '  +1.46%  ' },
    @{ Cell = 'D10'; Value = '0.0971' },
    @{ Cell = 'E10'; Value = '  -0.60%  ' },
    @{ Cell = 'E11'; Value = '  -1.00%  ' },
    @{ Cell = 'E12'; Value = '  -2.76%  ' },
    @{ Cell = 'E13'; Value = '  -2.22%  ' },
    @{ Cell = 'D14'; Value = '2.952.80' },
    @{ Cell = 'E14'; Value = '  +1.89%  ' },
    @{ Cell = 'D15'; Value = '58.253.08' },
    @{ Cell = 'E15'; Value = '  +0.44%  ' },
    @{ Cell = 'D16'; Value = '21.99' },
    @{ Cell = 'E16'; Value = '  -1.22%  ' },
    @{ Cell = 'D17'; Value = '0.0000134' },
    @{ Cell = 'E17'; Value = '  -0.28%  ' },
    @{ Cell = 'D18'; Value = '2.496.12' },
    @{ Cell = 'E18'; Value = '  +1.10%  ' },
    @{ Cell = 'D19'; Value = '10.58' },
    @{ Cell = 'E19'; Value = '  -0.23%  ' },
    @{ Cell = 'D20'; Value = '320.59' },
    @{ Cell = 'E20'; Value = '  +0.12%  ' },
    @{ Cell = 'D21'; Value = '4.14' },
    @{ Cell = 'E21'; Value = '  -0.13%  ' },
    @{ Cell = 'D22'; Value = '6.12' },
    @{ Cell = 'E22'; Value = '  +7.22%  ' },
    @{ Cell = 'D24'; Value = '64.30' },
    @{ Cell = 'E24'; Value = '  -0.54%  ' },
    @{ Cell = 'D26'; Value = '0.997' },
    @{ Cell = 'E26'; Value = '  -0.15%  ' },
    @{ Cell = 'E27'; Value = '  +0.34%  ' },
    @{ Cell = 'D28'; Value = '7.36' },
    @{ Cell = 'E28'; Value = '  +0.47%  ' },
    @{ Cell = 'D29'; Value = '0.0₃0750' },
    @{ Cell = 'E29'; Value = '  +0.54%  ' },
    @{ Cell = 'E30'; Value = '  +1.11%  ' },
    @{ Cell = 'D31'; Value = '167.15' },
    @{ Cell = 'E31'; Value = '  -0.10%  ' },
    @{ Cell = 'E32'; Value = '  +2.22%  ' },
    @{ Cell = 'D33'; Value = '6.25' },
    @{ Cell = 'E33'; Value = '  +0.65%  ' },
    @{ Cell = 'D34'; Value = '0.999' },
    @{ Cell = 'E34'; Value = '  +0.05%  ' },
    @{ Cell = 'D35'; Value = '0.998' },
    @{ Cell = 'E35'; Value = '  +0.11%  ' },
    @{ Cell = 'D36'; Value = '18.01' },
    @{ Cell = 'E36'; Value = '  +0.05%  ' },
    @{ Cell = 'E37'; Value = '  -8.54%  ' },
    @{ Cell = 'D38'; Value = '3.93' },
    @{ Cell = 'E38'; Value = '  -0.79%  ' },
    @{ Cell = 'D39'; Value = '1.46' },
    @{ Cell = 'E39'; Value = '  -0.09%  ' },
    @{ Cell = 'D40'; Value = '36.07' },
    @{ Cell = 'E40'; Value = '  -0.52%  ' },
    @{ Cell = 'D41'; Value = '0.769' },
    @{ Cell = 'E41'; Value = '  -3.09%  ' },
    @{ Cell = 'D42'; Value = '278.08' },
    @{ Cell = 'E42'; Value = '  +2.13%  ' },
    @{ Cell = 'D43'; Value = '3.47' },
    @{ Cell = 'E43'; Value = '  +0.86%  ' },
    @{ Cell = 'D44'; Value = '5.07' },
    @{ Cell = 'E44'; Value = '  +1.09%  ' },
    @{ Cell = 'D45'; Value = '0.593' },
    @{ Cell = 'E45'; Value = '  +0.73%  ' },
    @{ Cell = 'D46'; Value = '122.90' },
    @{ Cell = 'E46'; Value = '  -0.89%  ' },
    @{ Cell = 'D47'; Value = '0.0919' },
    @{ Cell = 'E47'; Value = '  +1.38%  ' },
    @{ Cell = 'E48'; Value = '  +2.54%  ' },
    @{ Cell = 'D49'; Value = '17.61' },
    @{ Cell = 'E49'; Value = '  +0.23%  ' },
    @{ Cell = 'D50'; Value = '0.0212' },
    @{ Cell = 'E50'; Value = '  +0.03%  ' },
    @{ Cell = 'D51'; Value = '16.70' },
    @{ Cell = 'E51'; Value = '  -0.78%  ' }
)

$helper = $ws.Range("H1")

foreach ($u in $updates) {
    $escaped = $u.Value -replace "'", "''"
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy()
    $ws.Range($u.Cell).PasteSpecial(-4163)
}

$helper.ClearContents()
$excel.CutCopyMode = 0
